$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "jumlah" column data (header text and value) while keeping the
# existing cell formatting/style on C1 and C2.
$ws.Range("C1").ClearContents()
$ws.Range("C2").ClearContents()

# Update the active selection to F2 (was B2).
$ws.Range("F2").Select()
